$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row: "<name>_old" -> "<name>_FV2310", "<name>_new" -> "<name>_FV2404"
$lastCol = 21
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value2
    if ($v -like "*_old") {
        $cell.Value = ($v -replace "_old$", "_FV2310")
    } elseif ($v -like "*_new") {
        $cell.Value = ($v -replace "_new$", "_FV2404")
    }
}

# 2. Freeze the header row (split after row 1, freeze top pane).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3. Turn the used range into a native Excel table ("Table1") covering the data.
$lastRow = $ws.UsedRange.Rows.Count
$tblRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $tblRange, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"
